# cleaned code and added comments
# Update column B ("Numerical Value") for rows 285-312 on Sheet1
# with refreshed sample values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> new value
$updates = @{
    285 = 1041
    286 = 4987
    287 = 6726
    288 = 1912
    289 = 5022
    290 = 7446
    291 = 8344
    292 = 6809
    293 = 9828
    294 = 7130
    295 = 2300
    296 = 9537
    297 = 8180
    298 = 6740
    299 = 5121
    300 = 2307
    301 = 7855
    302 = 3654
    303 = 8339
    304 = 2263
    305 = 2370
    306 = 1632
    307 = 3101
    308 = 4832
    309 = 7298
    310 = 9288
    311 = 1603
    312 = 1704
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
